$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.704.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.56'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4721'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3969'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.22'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08046'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.87'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.877.73'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.971'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.176'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.14'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001046'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06610'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.24'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.720.20'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.505'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.298'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.077.61'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.65'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.28'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.096'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.590'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.62'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9738'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09557'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.453'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.625'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.317'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06130'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02261'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.231'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.148'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6015'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1903'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.23'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5716'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.244'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.28'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.400'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.934'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000316'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06819'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.51%  '
